$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 85.5
$ws.Range("I9").Value = 47.714287
$ws.Range("J9").Value = 350
$ws.Range("K9").Value = 47.714287
$ws.Range("L9").Value = 350
$ws.Range("M9").Value = 121.285713
$ws.Range("N9").Value = -688

$ws.Range("H33").Value = 317.7037
$ws.Range("I33").Value = 257.45834
$ws.Range("K33").Value = 257.45834
$ws.Range("M33").Value = -28.45834000000002

$ws.Range("H113").Value = 55559416
$ws.Range("I113").Value = 125002980
$ws.Range("J113").Value = 4568.7
$ws.Range("K113").Value = 125002980
$ws.Range("L113").Value = 4568.7
$ws.Range("M113").Value = -124999726
$ws.Range("N113").Value = -11076.7

$ws.Range("H129").Value = 890.8197
$ws.Range("J129").Value = 897.6609999999999
$ws.Range("L129").Value = 2692.983
$ws.Range("N129").Value = -12692.983

$ws.Range("H136").Value = 49800
$ws.Range("J136").Value = 49800
$ws.Range("L136").Value = 49800
$ws.Range("N136").Value = -60000

$ws.Range("H137").Value = 1590.9166
$ws.Range("I137").Value = 1236.375
$ws.Range("J137").Value = 2300
$ws.Range("K137").Value = 3709.125
$ws.Range("L137").Value = 6900
$ws.Range("M137").Value = -1159.125
$ws.Range("N137").Value = -12000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3209
$ws.Range("I2").Value = 2475
$ws.Range("K2").Value = 2475
$ws.Range("M2").Value = -2362

$ws.Range("H116").Value = 3209
$ws.Range("I116").Value = 2475
$ws.Range("K116").Value = 2475
$ws.Range("M116").Value = -181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3209
$ws.Range("I3").Value = 2475
$ws.Range("K3").Value = 2475
$ws.Range("M3").Value = -2361

$ws.Range("H86").Value = 1450.25
$ws.Range("I86").Value = 1229.8
$ws.Range("J86").Value = 1817.6666
$ws.Range("K86").Value = 1229.8
$ws.Range("L86").Value = 1817.6666
$ws.Range("M86").Value = -106.8
$ws.Range("N86").Value = -4063.6666

$ws.Range("H89").Value = 1450.25
$ws.Range("I89").Value = 1229.8
$ws.Range("J89").Value = 1817.6666
$ws.Range("K89").Value = 6149
$ws.Range("L89").Value = 9088.333000000001
$ws.Range("M89").Value = -533
$ws.Range("N89").Value = -20320.333

$ws.Range("H94").Value = 5125.6665
$ws.Range("I94").Value = 1901.6
$ws.Range("J94").Value = 7428.5713
$ws.Range("K94").Value = 1901.6
$ws.Range("L94").Value = 7428.5713
$ws.Range("M94").Value = -1450.6
$ws.Range("N94").Value = -8330.5713

$ws.Range("H95").Value = 23324.6
$ws.Range("J95").Value = 23324.6
$ws.Range("L95").Value = 23324.6
$ws.Range("N95").Value = -28816.6

$ws.Range("H105").Value = 7144728.5
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H107").Value = 1405.125
$ws.Range("I107").Value = 1320.1428
$ws.Range("K107").Value = 1320.1428
$ws.Range("M107").Value = 599.8571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5447.647
$ws.Range("I99").Value = 4144.2856
$ws.Range("J99").Value = 6360
$ws.Range("K99").Value = 4144.2856
$ws.Range("L99").Value = 6360
$ws.Range("M99").Value = -2646.2856
$ws.Range("N99").Value = -9356

$ws.Range("H126").Value = 5447.647
$ws.Range("I126").Value = 4144.2856
$ws.Range("J126").Value = 6360
$ws.Range("K126").Value = 12432.8568
$ws.Range("L126").Value = 19080
$ws.Range("M126").Value = -9962.856800000001
$ws.Range("N126").Value = -24020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3809.3667
$ws.Range("I107").Value = 5164.2383
$ws.Range("J107").Value = 648
$ws.Range("K107").Value = 15492.7149
$ws.Range("L107").Value = 1944
$ws.Range("M107").Value = -13572.7149
$ws.Range("N107").Value = -5784

$ws.Range("H131").Value = 808.88
$ws.Range("I131").Value = 333.33334
$ws.Range("J131").Value = 823.5876500000001
$ws.Range("K131").Value = 1000.00002
$ws.Range("L131").Value = 2470.76295
$ws.Range("M131").Value = 4039.99998
$ws.Range("N131").Value = -12550.76295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 28745
$ws.Range("J39").Value = 28745
$ws.Range("L39").Value = 28745
$ws.Range("N39").Value = -29809

$ws.Range("H80").Value = 3079.375
$ws.Range("I80").Value = 2709.0908
$ws.Range("J80").Value = 3392.6924
$ws.Range("K80").Value = 2709.0908
$ws.Range("L80").Value = 3392.6924
$ws.Range("M80").Value = -1711.0908
$ws.Range("N80").Value = -5388.6924

$ws.Range("H83").Value = 3079.375
$ws.Range("I83").Value = 2709.0908
$ws.Range("J83").Value = 3392.6924
$ws.Range("K83").Value = 13545.454
$ws.Range("L83").Value = 16963.462
$ws.Range("M83").Value = -8553.454
$ws.Range("N83").Value = -26947.462

$ws.Range("H102").Value = 2243.75
$ws.Range("I102").Value = 2323.5557
$ws.Range("K102").Value = 2323.5557
$ws.Range("M102").Value = -701.5556999999999

$ws.Range("H113").Value = 2515.9473
$ws.Range("I113").Value = 2433.3333
$ws.Range("K113").Value = 2433.3333
$ws.Range("M113").Value = -263.3332999999998

$ws.Range("H126").Value = 4474.2856
$ws.Range("I126").Value = 3547.8262
$ws.Range("J126").Value = 6250
$ws.Range("K126").Value = 10643.4786
$ws.Range("L126").Value = 18750
$ws.Range("M126").Value = -8173.4786
$ws.Range("N126").Value = -23690

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8027.727
$ws.Range("I61").Value = 4575
$ws.Range("J61").Value = 10000.714
$ws.Range("K61").Value = 4575
$ws.Range("L61").Value = 10000.714
$ws.Range("M61").Value = -4373
$ws.Range("N61").Value = -10404.714

$ws.Range("H82").Value = 2244.9
$ws.Range("I82").Value = 3490
$ws.Range("J82").Value = 999.8
$ws.Range("K82").Value = 3490
$ws.Range("L82").Value = 999.8
$ws.Range("M82").Value = -3129
$ws.Range("N82").Value = -1721.8

$ws.Range("H85").Value = 2244.9
$ws.Range("I85").Value = 3490
$ws.Range("J85").Value = 999.8
$ws.Range("K85").Value = 3490
$ws.Range("L85").Value = 999.8
$ws.Range("M85").Value = -2242
$ws.Range("N85").Value = -3495.8

$ws.Range("H113").Value = 8027.727
$ws.Range("I113").Value = 4575
$ws.Range("J113").Value = 10000.714
$ws.Range("K113").Value = 4575
$ws.Range("L113").Value = 10000.714
$ws.Range("M113").Value = -2405
$ws.Range("N113").Value = -14340.714

$ws.Range("H140").Value = 46633.5
$ws.Range("J140").Value = 46633.5
$ws.Range("L140").Value = 46633.5
$ws.Range("N140").Value = -56993.5
